# Update absenteeism data rows 2-11 per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = 58378
$ws.Range("B2").Value = "Evelyn Nogueira"
$ws.Range("C2").Value = "Recursos Humanos"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45094
$ws.Range("G2").Value = 8675.889999999999

# Row 3
$ws.Range("A3").Value = 23356
$ws.Range("B3").Value = "Camila Silveira"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45090
$ws.Range("G3").Value = 9222.18

# Row 4
$ws.Range("A4").Value = 44933
$ws.Range("B4").Value = "Maria Alice Lopes"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45091
$ws.Range("G4").Value = 5703.94

# Row 5
$ws.Range("A5").Value = 14556
$ws.Range("B5").Value = "Noah Cardoso"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45105
$ws.Range("G5").Value = 10882.03

# Row 6
$ws.Range("A6").Value = 8025
$ws.Range("B6").Value = "Thiago Barros"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45095
$ws.Range("G6").Value = 12436.01

# Row 7
$ws.Range("A7").Value = 42418
$ws.Range("B7").Value = "Thiago Viana"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45101
$ws.Range("G7").Value = 5203.48

# Row 8
$ws.Range("A8").Value = 80293
$ws.Range("B8").Value = "Yasmin Costa"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45079
$ws.Range("G8").Value = 7933.98

# Row 9
$ws.Range("A9").Value = 49639
$ws.Range("B9").Value = "Luiz Felipe Campos"
$ws.Range("F9").Value = 45099
$ws.Range("G9").Value = 4297.17

# Row 10
$ws.Range("A10").Value = 97251
$ws.Range("B10").Value = "Evelyn da Paz"
$ws.Range("C10").Value = "Marketing"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 6403.86

# Row 11
$ws.Range("A11").Value = 1761
$ws.Range("B11").Value = "Gustavo Cardoso"
$ws.Range("D11").Value = "Outros"
$ws.Range("F11").Value = 45078
$ws.Range("G11").Value = 7977.81
